$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: E1 -> pipeline_id, F1 -> config_json (new), G1 -> updated_at (new)
$ws.Range("E1").Value = "pipeline_id"
$ws.Range("F1").Value = "config_json"
$ws.Range("G1").Value = "updated_at"

# Copy header style (bold, centered, bordered) from F1 (old style) to the new G1 cell
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Remove the old data row (row 2), leaving just the header row
$ws.Rows("2:2").Delete() | Out-Null
